$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Median_Pulse_Width (column C) values for the Mounted Tire Processing Pipeline
    $ws.Cells.Item(6, 3).Value = 12
    $ws.Cells.Item(7, 3).Value = 19
    $ws.Cells.Item(8, 3).Value = 21
    $ws.Cells.Item(9, 3).Value = 23
    $ws.Cells.Item(31, 3).Value = 16
    $ws.Cells.Item(32, 3).Value = 17
    $ws.Cells.Item(33, 3).Value = 19
    $ws.Cells.Item(39, 3).Value = 10
    $ws.Cells.Item(40, 3).Value = 12
    $ws.Cells.Item(41, 3).Value = 21
    $ws.Cells.Item(42, 3).Value = 11
    $ws.Cells.Item(43, 3).Value = 15
    $ws.Cells.Item(44, 3).Value = 17
    $ws.Cells.Item(46, 3).Value = 9
    $ws.Cells.Item(47, 3).Value = 17
    $ws.Cells.Item(50, 3).Value = 8
    $ws.Cells.Item(52, 3).Value = 18
    $ws.Cells.Item(53, 3).Value = 25
    $ws.Cells.Item(56, 3).Value = 18
    $ws.Cells.Item(57, 3).Value = 22
    $ws.Cells.Item(66, 3).Value = 5
    $ws.Cells.Item(67, 3).Value = 13
    $ws.Cells.Item(68, 3).Value = 18
    $ws.Cells.Item(69, 3).Value = 25
    $ws.Cells.Item(82, 3).Value = 10
    $ws.Cells.Item(84, 3).Value = 16
    $ws.Cells.Item(85, 3).Value = 18
    $ws.Cells.Item(86, 3).Value = 8
    $ws.Cells.Item(87, 3).Value = 12
    $ws.Cells.Item(92, 3).Value = 19
    $ws.Cells.Item(93, 3).Value = 25
    $ws.Cells.Item(94, 3).Value = 6
    $ws.Cells.Item(95, 3).Value = 14
    $ws.Cells.Item(96, 3).Value = 14
    $ws.Cells.Item(97, 3).Value = 16
    $ws.Cells.Item(100, 3).Value = 21
    $ws.Cells.Item(101, 3).Value = 26
    $ws.Cells.Item(104, 3).Value = 21
    $ws.Cells.Item(105, 3).Value = 26
    $ws.Cells.Item(110, 3).Value = 6
    $ws.Cells.Item(111, 3).Value = 13
    $ws.Cells.Item(112, 3).Value = 14
    $ws.Cells.Item(113, 3).Value = 18
    $ws.Cells.Item(114, 3).Value = 8
    $ws.Cells.Item(134, 3).Value = 12
    $ws.Cells.Item(135, 3).Value = 17
    $ws.Cells.Item(136, 3).Value = 19
    $ws.Cells.Item(137, 3).Value = 26
    $ws.Cells.Item(169, 3).Value = 25
    $ws.Cells.Item(196, 3).Value = 14
    $ws.Cells.Item(197, 3).Value = 22
    $ws.Cells.Item(216, 3).Value = 12
    $ws.Cells.Item(218, 3).Value = 7
    $ws.Cells.Item(219, 3).Value = 12
    $ws.Cells.Item(220, 3).Value = 13
    $ws.Cells.Item(221, 3).Value = 17
    $ws.Cells.Item(222, 3).Value = 11
    $ws.Cells.Item(223, 3).Value = 17
    $ws.Cells.Item(224, 3).Value = 22
    $ws.Cells.Item(225, 3).Value = 28
    $ws.Cells.Item(226, 3).Value = 11
    $ws.Cells.Item(227, 3).Value = 16
    $ws.Cells.Item(228, 3).Value = 19
    $ws.Cells.Item(229, 3).Value = 22
